# changes done 1 Feb 2024
#
# The 30 Jan 2024 run's generated Policy/Report export names on the
# "policyNumber" sheet are bumped to the 1 Feb 2024 run, and the
# "policyNumber" sheet (rather than "loginCredentials") is left as the
# active/selected sheet & cell when the workbook is saved.

$wb = $excel.ActiveWorkbook
$wsPolicy = $wb.Worksheets.Item("policyNumber")

# B2/B5 hold the generated Policy_/Report_ file-name stems, timestamped
# with the run date (ddMMyyyy) - move them from 30012024 to 01022024.
$wsPolicy.Range("B2").Value = "Policy_01022024_"
$wsPolicy.Range("B5").Value = "Report_01022024_"

# Leave the selection on policyNumber!B10 and make that sheet the active
# (displayed) tab, rather than loginCredentials.
$wsPolicy.Range("B10").Select()
$wsPolicy.Activate()
